# Auto-generated Excel COM-interop script applying the Ultros_Profits diff.
# Updates LeveProfitNQ/LeveProfitHQ-related computed columns (H-N) across
# the ALC, ARM, BSM, CRP, CUL, LTW, WVR leve-profit tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 44.666668
$ws.Range("J5").Value = 44
$ws.Range("L5").Value = 44
$ws.Range("N5").Value = -274

$ws.Range("H18").Value = 4914.2856
$ws.Range("J18").Value = 7000
$ws.Range("L18").Value = 7000
$ws.Range("N18").Value = -7568

$ws.Range("H40").Value = 4128.7144

$ws.Range("H86").Value = 3257.5
$ws.Range("J86").Value = 3877.6667
$ws.Range("L86").Value = 3877.6667
$ws.Range("N86").Value = -6123.6667

$ws.Range("H89").Value = 3257.5
$ws.Range("J89").Value = 3877.6667
$ws.Range("L89").Value = 19388.3335
$ws.Range("N89").Value = -30620.3335

$ws.Range("H116").Value = 6508.115
$ws.Range("I116").Value = 5309.6665
$ws.Range("J116").Value = 7535.357
$ws.Range("K116").Value = 5309.6665
$ws.Range("L116").Value = 7535.357
$ws.Range("M116").Value = -1867.6665
$ws.Range("N116").Value = -14419.357

$ws.Range("H125").Value = 1879.4
$ws.Range("I125").Value = 1879.4
$ws.Range("K125").Value = 16914.6
$ws.Range("M125").Value = -14454.6

$ws.Range("H130").Value = 19948.125
$ws.Range("J130").Value = 19948.125
$ws.Range("L130").Value = 19948.125
$ws.Range("N130").Value = -29988.125

$ws.Range("H138").Value = 3167.8438
$ws.Range("I138").Value = 1529
$ws.Range("J138").Value = 4289.1577
$ws.Range("K138").Value = 4587
$ws.Range("L138").Value = 12867.4731
$ws.Range("M138").Value = 553
$ws.Range("N138").Value = -23147.4731

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4302.522
$ws.Range("I61").Value = 2560.9
$ws.Range("K61").Value = 2560.9
$ws.Range("M61").Value = -2348.9

$ws.Range("H132").Value = 27029384
$ws.Range("I132").Value = 34484988
$ws.Range("K132").Value = 103454964
$ws.Range("M132").Value = -103452434

$ws.Range("H136").Value = 4302.522
$ws.Range("I136").Value = 2560.9
$ws.Range("K136").Value = 7682.700000000001
$ws.Range("M136").Value = -5132.700000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2570
$ws.Range("I134").Value = 2570
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 7710
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -5175
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 402.1875
$ws.Range("I22").Value = 437.25
$ws.Range("J22").Value = 367.125
$ws.Range("K22").Value = 437.25
$ws.Range("L22").Value = 367.125
$ws.Range("M22").Value = -87.25
$ws.Range("N22").Value = -1067.125

$ws.Range("H48").Value = 46600
$ws.Range("J48").Value = 46600
$ws.Range("L48").Value = 46600
$ws.Range("N48").Value = -47552

$ws.Range("H86").Value = 56121
$ws.Range("J86").Value = 21281.285
$ws.Range("L86").Value = 21281.285
$ws.Range("N86").Value = -23527.285

$ws.Range("H89").Value = 56121
$ws.Range("J89").Value = 21281.285
$ws.Range("L89").Value = 106406.425
$ws.Range("N89").Value = -117638.425

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 235.07692
$ws.Range("I7").Value = 323.42856
$ws.Range("K7").Value = 970.28568
$ws.Range("M7").Value = -858.28568

$ws.Range("H11").Value = 84334.73
$ws.Range("I11").Value = 115876.5
$ws.Range("K11").Value = 347629.5
$ws.Range("M11").Value = -347489.5

$ws.Range("H25").Value = 1284.6154
$ws.Range("I25").Value = 962.5
$ws.Range("J25").Value = 1800
$ws.Range("K25").Value = 2887.5
$ws.Range("L25").Value = 5400
$ws.Range("M25").Value = -2718.5
$ws.Range("N25").Value = -5738

$ws.Range("H30").Value = 1284.6154
$ws.Range("I30").Value = 962.5
$ws.Range("J30").Value = 1800
$ws.Range("K30").Value = 2887.5
$ws.Range("L30").Value = 5400
$ws.Range("M30").Value = -2785.5
$ws.Range("N30").Value = -5604

$ws.Range("H107").Value = 10755106
$ws.Range("I107").Value = 2372.0908
$ws.Range("J107").Value = 16669110
$ws.Range("K107").Value = 7116.2724
$ws.Range("L107").Value = 50007330
$ws.Range("M107").Value = -5196.2724
$ws.Range("N107").Value = -50011170

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 7576946.5
$ws.Range("I22").Value = 11364424
$ws.Range("K22").Value = 11364424
$ws.Range("M22").Value = -11364129

$ws.Range("H27").Value = 7576946.5
$ws.Range("I27").Value = 11364424
$ws.Range("K27").Value = 11364424
$ws.Range("M27").Value = -11364317

$ws.Range("H46").Value = 2381.5789
$ws.Range("I46").Value = 1845.4286
$ws.Range("K46").Value = 1845.4286
$ws.Range("M46").Value = -1657.4286

$ws.Range("H61").Value = 999
$ws.Range("I61").Value = 1018.4762
$ws.Range("K61").Value = 1018.4762
$ws.Range("M61").Value = -816.4761999999999

$ws.Range("H82").Value = 58825540
$ws.Range("I82").Value = 83335416
$ws.Range("J82").Value = 1838.6
$ws.Range("K82").Value = 83335416
$ws.Range("L82").Value = 1838.6
$ws.Range("M82").Value = -83335055
$ws.Range("N82").Value = -2560.6

$ws.Range("H85").Value = 58825540
$ws.Range("I85").Value = 83335416
$ws.Range("J85").Value = 1838.6
$ws.Range("K85").Value = 83335416
$ws.Range("L85").Value = 1838.6
$ws.Range("M85").Value = -83334168
$ws.Range("N85").Value = -4334.6

$ws.Range("H93").Value = 1514
$ws.Range("I93").Value = 1123.3334
$ws.Range("J93").Value = 2100
$ws.Range("K93").Value = 1123.3334
$ws.Range("L93").Value = 2100
$ws.Range("M93").Value = 124.6666
$ws.Range("N93").Value = -4596

$ws.Range("H100").Value = 127901.11
$ws.Range("I100").Value = 224022
$ws.Range("K100").Value = 224022
$ws.Range("M100").Value = -223481

$ws.Range("H113").Value = 999
$ws.Range("I113").Value = 1018.4762
$ws.Range("K113").Value = 1018.4762
$ws.Range("M113").Value = 1151.5238

$ws.Range("H132").Value = 3044
$ws.Range("I132").Value = 2469.111
$ws.Range("K132").Value = 7407.333
$ws.Range("M132").Value = -4877.333

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 99997.5
$ws.Range("J92").Value = 99997.5
$ws.Range("L92").Value = 99997.5
$ws.Range("N92").Value = -104989.5

$ws.Range("H98").Value = 52293.75
$ws.Range("J98").Value = 52293.75
$ws.Range("L98").Value = 52293.75
$ws.Range("N98").Value = -58283.75

$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

$ws.Range("H132").Value = 2933.5454
$ws.Range("I132").Value = 2601.8696
$ws.Range("J132").Value = 3696.4
$ws.Range("K132").Value = 7805.6088
$ws.Range("L132").Value = 11089.2
$ws.Range("M132").Value = -5275.6088
$ws.Range("N132").Value = -16149.2
